# pridana tabulka subjective ratings
# Fill in the subjective-ratings scale labels / re-worded statements for
# rows 2-10 (columns B..E) and tidy up the view (active cell + column width).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "Rozhranie sa pouzivalo jednoducho" (unchanged wording, restated)
$ws.Range("B2").Value = "Rozhranie sa použivalo jednoducho"
$ws.Range("C2").Value = "veľmi ľahké"
$ws.Range("D2").Value = "6 5 4 3 2 1 0 "
$ws.Range("E2").Value = "veľmi ťažké"

# Row 3 - "Rozhranie bolo lahko pochopitelne" (unchanged wording, restated)
$ws.Range("B3").Value = "Rozhranie bolo ľahko pochopiteľne"
$ws.Range("C3").Value = "veľmi ľahké"
$ws.Range("D3").Value = "6 5 4 3 2 1 0 "
$ws.Range("E3").Value = "veľmi ťažké"

# Row 4 - new scale label "velmi jednoduche" for the left anchor
$ws.Range("B4").Value = "Nebolo jednoduché zistiť ako postupovať"
$ws.Range("C4").Value = "veľmi jednoduché"
$ws.Range("D4").Value = "6 5 4 3 2 1 0 "
$ws.Range("E4").Value = "veľmi ťažké"

# Row 5 - new scale labels "velmi jasne" / "velmi nejasne"
$ws.Range("B5").Value = "Nebolo jasné čo je treba urobiť"
$ws.Range("C5").Value = "veľmi jasné"
$ws.Range("D5").Value = "6 5 4 3 2 1 0 "
$ws.Range("E5").Value = "veľmi nejasné"

# Row 6
$ws.Range("B6").Value = "Na 3 stranke nebolo jasné ako postupovať ďalej"
$ws.Range("C6").Value = "veľmi jasné"
$ws.Range("D6").Value = "6 5 4 3 2 1 0 "
$ws.Range("E6").Value = "veľmi nejasné"

# Row 7
$ws.Range("B7").Value = "Na 2 stranke nebolo jasné ako postupovať ďalej"
$ws.Range("C7").Value = "veľmi jasné"
$ws.Range("D7").Value = "6 5 4 3 2 1 0 "
$ws.Range("E7").Value = "veľmi nejasné"

# Row 8
$ws.Range("B8").Value = "Na 1 stranke nebolo jasné ako postupovať ďalej"
$ws.Range("C8").Value = "veľmi jasné"
$ws.Range("D8").Value = "6 5 4 3 2 1 0 "
$ws.Range("E8").Value = "veľmi nejasné"

# Row 9 - new scale labels "velmi" / "vobec"
$ws.Range("B9").Value = "Chcel/a by som využívať toto rozhranie aj v budúcnosti"
$ws.Range("C9").Value = "veľmi"
$ws.Range("D9").Value = "6 5 4 3 2 1 0 "
$ws.Range("E9").Value = "vôbec"

# Row 10 - new scale labels "ziadne zmeny" / "velke zmeny"
$ws.Range("B10").Value = "Urobil/a by som veľké zmeny v tomto rozhraní"
$ws.Range("C10").Value = "žiadne zmeny"
$ws.Range("D10").Value = "6 5 4 3 2 1 0 "
$ws.Range("E10").Value = "veľké zmeny"

# Widen column C slightly to fit the new labels
$ws.Range("C1").ColumnWidth = 14.6

# Leave the cursor on E8, matching where editing finished
$ws.Range("E8").Select()
